$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as plain text so
# dot-separated numbers (e.g. "28.653.60") and other price strings
# are preserved exactly rather than being reinterpreted as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.653.60"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.952.34"
$ws.Range("E3").Value = "  -2.32%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.015"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("E5").Value = "  -2.60%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.014"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4771"
$ws.Range("E7").Value = "  -4.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4044"
$ws.Range("E8").Value = "  -4.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.76"
$ws.Range("E9").Value = "  +0.24%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08462"
$ws.Range("E10").Value = "  -5.87%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.063"
$ws.Range("E11").Value = "  -4.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.22"
$ws.Range("E12").Value = "  -4.71%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.960.24"
$ws.Range("E13").Value = "  -3.00%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.625"
$ws.Range("E14").Value = "  -5.22%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.202"
$ws.Range("E15").Value = "  -4.02%  "

$ws.Range("E16").Value = "  +0.10%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "89.38"
$ws.Range("E17").Value = "  -4.73%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001072"
$ws.Range("E18").Value = "  -3.46%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06628"
$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.65"
$ws.Range("E20").Value = "  -5.31%  "

$ws.Range("E21").Value = "  +0.02%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.820"
$ws.Range("E22").Value = "  -2.25%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.675.43"
$ws.Range("E23").Value = "  -3.11%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.53"
$ws.Range("E24").Value = "  -3.59%  "

$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.188.82"
$ws.Range("E26").Value = "  -3.21%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.93"
$ws.Range("E27").Value = "  -2.77%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").Value = "  -2.26%  "

$ws.Range("E29").Value = "  -7.57%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.155"
$ws.Range("E30").Value = "  -6.12%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.15"
$ws.Range("E31").Value = "  -3.13%  "

$ws.Range("E32").Value = "  -4.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09584"
$ws.Range("E33").Value = "  -3.53%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.659"
$ws.Range("E34").Value = "  -2.83%  "

$ws.Range("E35").Value = "  -3.49%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.435"
$ws.Range("E36").Value = "  -8.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02331"
$ws.Range("E37").Value = "  -5.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06229"
$ws.Range("E38").Value = "  -1.78%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.267"
$ws.Range("E39").Value = "  -3.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.713"
$ws.Range("E40").Value = "  -6.81%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6237"
$ws.Range("E41").Value = "  -4.80%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.09"
$ws.Range("E42").Value = "  -4.88%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.012"
$ws.Range("E43").Value = "  +0.04%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1922"
$ws.Range("E44").Value = "  -6.18%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.341"
$ws.Range("E45").Value = "  +2.95%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5943"
$ws.Range("E46").Value = "  -6.04%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.97"
$ws.Range("E47").Value = "  -3.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.084"
$ws.Range("E48").Value = "  -4.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.405"
$ws.Range("E49").Value = "  -3.12%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000335"
$ws.Range("E50").Value = "  -0.55%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06822"
$ws.Range("E51").Value = "  -2.34%  "
